$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "K8JqOFWo"
$ws.Cells.Item(2, 3).Value = "09:30"
$ws.Cells.Item(2, 4).Value = "MALAYSIA - SUPER LEAGUE"
$ws.Cells.Item(2, 5).Value = "Negeri Sembilan"
$ws.Cells.Item(2, 6).Value = "Johor DT"
$ws.Cells.Item(2, 7).Value = 24
$ws.Cells.Item(2, 8).Value = 9.5
$ws.Cells.Item(2, 9).Value = 1.04
$ws.Cells.Item(2, 10).Value = 18
$ws.Cells.Item(2, 11).Value = 3.85
$ws.Cells.Item(2, 12).Value = 1.23
$ws.Cells.Item(2, 13).Value = 60
$ws.Cells.Item(2, 14).Value = 26
$ws.Cells.Item(2, 15).Value = 1.01
$ws.Cells.Item(2, 16).Value = 11
$ws.Cells.Item(2, 17).Value = 1.16
$ws.Cells.Item(2, 18).Value = 3.88
$ws.Cells.Item(2, 19).Value = 1.09
$ws.Cells.Item(2, 20).Value = 6.1
$ws.Cells.Item(2, 21).Value = 2.84
$ws.Cells.Item(2, 22).Value = 1.4
$ws.Cells.Item(2, 23).Value = 90
$ws.Cells.Item(2, 24).Value = 400
$ws.Cells.Item(2, 25).Value = 120
$ws.Cells.Item(2, 26).Value = 500
$ws.Cells.Item(2, 27).Value = 400
$ws.Cells.Item(2, 28).Value = 400
$ws.Cells.Item(2, 29).Value = 27
$ws.Cells.Item(2, 30).Value = 28
$ws.Cells.Item(2, 31).Value = 65
$ws.Cells.Item(2, 32).Value = 300
$ws.Cells.Item(2, 33).Value = 201
$ws.Cells.Item(2, 34).Value = 12.5
$ws.Cells.Item(2, 35).Value = 6.9
$ws.Cells.Item(2, 36).Value = 15.5
$ws.Cells.Item(2, 37).Value = 5.4
$ws.Cells.Item(2, 38).Value = 12
$ws.Cells.Item(2, 39).Value = 50
$ws.Cells.Item(2, 40).Value = 27
$ws.Cells.Item(2, 41).Value = 250
$ws.Cells.Item(2, 42).Value = 150
$ws.Cells.Item(2, 43).Value = 501
$ws.Cells.Item(2, 44).Value = 501
$ws.Cells.Item(2, 45).Value = 501
$ws.Cells.Item(2, 46).Value = 5.4
$ws.Cells.Item(2, 47).Value = 16
$ws.Cells.Item(2, 48).Value = 175
$ws.Cells.Item(2, 49).Value = 3.3
$ws.Cells.Item(2, 50).Value = 3.65
$ws.Cells.Item(2, 51).Value = 18
$ws.Cells.Item(2, 52).Value = 5.9
$ws.Cells.Item(2, 53).Value = 29
$ws.Cells.Item(2, 54).Value = 250
$ws.Cells.Item(2, 55).Value = 51
$ws.Cells.Item(2, 56).Value = 51

# Row 3
$ws.Cells.Item(3, 1).Value = "SYMiMg1b"
$ws.Cells.Item(3, 3).Value = "10:00"
$ws.Cells.Item(3, 4).Value = "MALAYSIA - SUPER LEAGUE"
$ws.Cells.Item(3, 5).Value = "Terengganu"
$ws.Cells.Item(3, 6).Value = "Penang"
$ws.Cells.Item(3, 7).Value = 1.31
$ws.Cells.Item(3, 8).Value = 4.4
$ws.Cells.Item(3, 9).Value = 8.25
$ws.Cells.Item(3, 10).Value = 1.78
$ws.Cells.Item(3, 11).Value = 2.4
$ws.Cells.Item(3, 12).Value = 6.8
$ws.Cells.Item(3, 13).Value = 1.02
$ws.Cells.Item(3, 14).Value = 13.8
$ws.Cells.Item(3, 15).Value = 1.15
$ws.Cells.Item(3, 16).Value = 4
$ws.Cells.Item(3, 17).Value = 1.53
$ws.Cells.Item(3, 18).Value = 2.18
$ws.Cells.Item(3, 19).Value = 1.3
$ws.Cells.Item(3, 20).Value = 3.32
$ws.Cells.Item(3, 21).Value = 1.84
$ws.Cells.Item(3, 22).Value = 1.92
$ws.Cells.Item(3, 23).Value = 6.5
$ws.Cells.Item(3, 24).Value = 5.8
$ws.Cells.Item(3, 25).Value = 7
$ws.Cells.Item(3, 26).Value = 7.3
$ws.Cells.Item(3, 27).Value = 8.75
$ws.Cells.Item(3, 28).Value = 19
$ws.Cells.Item(3, 29).Value = 13.5
$ws.Cells.Item(3, 30).Value = 8
$ws.Cells.Item(3, 31).Value = 15
$ws.Cells.Item(3, 32).Value = 55
$ws.Cells.Item(3, 33).Value = 300
$ws.Cells.Item(3, 34).Value = 21
$ws.Cells.Item(3, 35).Value = 50
$ws.Cells.Item(3, 36).Value = 20
$ws.Cells.Item(3, 37).Value = 175
$ws.Cells.Item(3, 38).Value = 70
$ws.Cells.Item(3, 39).Value = 50
$ws.Cells.Item(3, 40).Value = 3.2
$ws.Cells.Item(3, 41).Value = 5.9
$ws.Cells.Item(3, 42).Value = 14.5
$ws.Cells.Item(3, 43).Value = 15.5
$ws.Cells.Item(3, 44).Value = 40
$ws.Cells.Item(3, 45).Value = 175
$ws.Cells.Item(3, 46).Value = 3.2
$ws.Cells.Item(3, 47).Value = 7.8
$ws.Cells.Item(3, 48).Value = 65
$ws.Cells.Item(3, 49).Value = 9.25
$ws.Cells.Item(3, 50).Value = 45
$ws.Cells.Item(3, 51).Value = 40
$ws.Cells.Item(3, 52).Value = 300
$ws.Cells.Item(3, 53).Value = 300
$ws.Cells.Item(3, 54).Value = 450
$ws.Cells.Item(3, 55).Value = 51
$ws.Cells.Item(3, 56).Value = 51

# Row 4
$ws.Cells.Item(4, 1).Value = "8QJqec8E"
$ws.Cells.Item(4, 3).Value = "11:30"
$ws.Cells.Item(4, 4).Value = "QATAR - QSL"
$ws.Cells.Item(4, 5).Value = "Shamal"
$ws.Cells.Item(4, 6).Value = "Al-Duhail"
$ws.Cells.Item(4, 7).Value = 5.6
$ws.Cells.Item(4, 8).Value = 4.7
$ws.Cells.Item(4, 9).Value = 1.45
$ws.Cells.Item(4, 10).Value = 5.1
$ws.Cells.Item(4, 11).Value = 2.55
$ws.Cells.Item(4, 12).Value = 1.91
$ws.Cells.Item(4, 13).Value = 1.02
$ws.Cells.Item(4, 14).Value = 9.75
$ws.Cells.Item(4, 15).Value = 1.14
$ws.Cells.Item(4, 16).Value = 5
$ws.Cells.Item(4, 17).Value = 1.44
$ws.Cells.Item(4, 18).Value = 2.62
$ws.Cells.Item(4, 19).Value = 1.26
$ws.Cells.Item(4, 20).Value = 3.55
$ws.Cells.Item(4, 21).Value = 1.57
$ws.Cells.Item(4, 22).Value = 2.25
$ws.Cells.Item(4, 23).Value = 22
$ws.Cells.Item(4, 24).Value = 40
$ws.Cells.Item(4, 25).Value = 17.5
$ws.Cells.Item(4, 26).Value = 100
$ws.Cells.Item(4, 27).Value = 45
$ws.Cells.Item(4, 28).Value = 40
$ws.Cells.Item(4, 29).Value = 9.75
$ws.Cells.Item(4, 30).Value = 9.75
$ws.Cells.Item(4, 31).Value = 15
$ws.Cells.Item(4, 32).Value = 50
$ws.Cells.Item(4, 33).Value = 250
$ws.Cells.Item(4, 34).Value = 10.5
$ws.Cells.Item(4, 35).Value = 9
$ws.Cells.Item(4, 36).Value = 8.5
$ws.Cells.Item(4, 37).Value = 11.25
$ws.Cells.Item(4, 38).Value = 10.5
$ws.Cells.Item(4, 39).Value = 18.5
$ws.Cells.Item(4, 40).Value = 7.4
$ws.Cells.Item(4, 41).Value = 28
$ws.Cells.Item(4, 42).Value = 28
$ws.Cells.Item(4, 43).Value = 150
$ws.Cells.Item(4, 44).Value = 150
$ws.Cells.Item(4, 45).Value = 300
$ws.Cells.Item(4, 46).Value = 3.55
$ws.Cells.Item(4, 47).Value = 7.3
$ws.Cells.Item(4, 48).Value = 50
$ws.Cells.Item(4, 49).Value = 3.65
$ws.Cells.Item(4, 50).Value = 6.7
$ws.Cells.Item(4, 51).Value = 13.5
$ws.Cells.Item(4, 52).Value = 17.5
$ws.Cells.Item(4, 53).Value = 35
$ws.Cells.Item(4, 54).Value = 150
$ws.Cells.Item(4, 55).Value = 51
$ws.Cells.Item(4, 56).Value = 51

# Row 5
$ws.Cells.Item(5, 1).Value = "QXl1xs0B"
$ws.Cells.Item(5, 3).Value = "11:35"
$ws.Cells.Item(5, 4).Value = "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE"
$ws.Cells.Item(5, 5).Value = "Al Khaleej"
$ws.Cells.Item(5, 6).Value = "Al Raed"
$ws.Cells.Item(5, 7).Value = 2.15
$ws.Cells.Item(5, 8).Value = 3.3
$ws.Cells.Item(5, 9).Value = 3.3
$ws.Cells.Item(5, 10).Value = 2.75
$ws.Cells.Item(5, 11).Value = 2.05
$ws.Cells.Item(5, 12).Value = 3.75
$ws.Cells.Item(5, 13).Value = 1.06
$ws.Cells.Item(5, 14).Value = 8
$ws.Cells.Item(5, 15).Value = 1.33
$ws.Cells.Item(5, 16).Value = 3.25
$ws.Cells.Item(5, 17).Value = 2.05
$ws.Cells.Item(5, 18).Value = 1.75
$ws.Cells.Item(5, 19).Value = 1.44
$ws.Cells.Item(5, 20).Value = 2.63
$ws.Cells.Item(5, 21).Value = 1.83
$ws.Cells.Item(5, 22).Value = 1.83
$ws.Cells.Item(5, 23).Value = 7.5
$ws.Cells.Item(5, 24).Value = 10
$ws.Cells.Item(5, 25).Value = 9.5
$ws.Cells.Item(5, 26).Value = 19
$ws.Cells.Item(5, 27).Value = 19
$ws.Cells.Item(5, 28).Value = 29
$ws.Cells.Item(5, 29).Value = 9
$ws.Cells.Item(5, 30).Value = 6.5
$ws.Cells.Item(5, 31).Value = 15
$ws.Cells.Item(5, 32).Value = 51
$ws.Cells.Item(5, 33).Value = 700
$ws.Cells.Item(5, 34).Value = 9.5
$ws.Cells.Item(5, 35).Value = 17
$ws.Cells.Item(5, 36).Value = 12
$ws.Cells.Item(5, 37).Value = 34
$ws.Cells.Item(5, 38).Value = 29
$ws.Cells.Item(5, 39).Value = 41
$ws.Cells.Item(5, 40).Value = 4.33
$ws.Cells.Item(5, 41).Value = 12
$ws.Cells.Item(5, 42).Value = 23
$ws.Cells.Item(5, 43).Value = 41
$ws.Cells.Item(5, 44).Value = 67
$ws.Cells.Item(5, 45).Value = 151
$ws.Cells.Item(5, 46).Value = 2.63
$ws.Cells.Item(5, 47).Value = 8
$ws.Cells.Item(5, 48).Value = 51
$ws.Cells.Item(5, 49).Value = 5
$ws.Cells.Item(5, 50).Value = 19
$ws.Cells.Item(5, 51).Value = 29
$ws.Cells.Item(5, 52).Value = 51
$ws.Cells.Item(5, 53).Value = 81
$ws.Cells.Item(5, 54).Value = 300
$ws.Cells.Item(5, 55).Value = 81
$ws.Cells.Item(5, 56).Value = 81

# Row 7 odds tweak
$ws.Cells.Item(7, 17).Value = 2.1
$ws.Cells.Item(7, 18).Value = 1.7
